# Daily Report update: 2026-01-06
# Appends a new day's block (date serial 46027) to the Daily_Data sheet,
# mirroring the existing 22-row block (rows 2-23) into rows 24-45.
# For each institution/region row, the new day's PREV_TOTAL (C) and
# TOTAL_TODAY (H) carry forward the prior day's TOTAL_TODAY (H), while
# RECEIVED/WITHDRAWN/NET_CHANGE/ADJUSTMENT (D:G) reset to 0 for the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

$newDateSerial = 46027
$firstDataRow  = 2
$lastDataRow   = 23
$rowOffset     = $lastDataRow - $firstDataRow + 1   # 22 rows per day block
$dateNumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $newRow = $row + $rowOffset

    $regionType  = $ws.Cells.Item($row, 2).Value2   # B: Region_Type
    $priorTotal  = $ws.Cells.Item($row, 8).Value2   # H: prior day's TOTAL_TODAY

    # A: Date
    $ws.Cells.Item($newRow, 1).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($newRow, 1).Value2 = $newDateSerial
    # B: Region_Type
    $ws.Cells.Item($newRow, 2).Value2 = $regionType
    # C: PREV_TOTAL = prior day's TOTAL_TODAY
    $ws.Cells.Item($newRow, 3).Value2 = $priorTotal
    # D: RECEIVED
    $ws.Cells.Item($newRow, 4).Value2 = 0
    # E: WITHDRAWN
    $ws.Cells.Item($newRow, 5).Value2 = 0
    # F: NET_CHANGE
    $ws.Cells.Item($newRow, 6).Value2 = 0
    # G: ADJUSTMENT
    $ws.Cells.Item($newRow, 7).Value2 = 0
    # H: TOTAL_TODAY = PREV_TOTAL + NET_CHANGE + ADJUSTMENT (no activity today)
    $ws.Cells.Item($newRow, 8).Value2 = $priorTotal
}
